$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.146.42"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "2.226.39"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'291.76"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").Value = "'87.84"
$ws.Range("E6").Value = "  +1.75%  "

$ws.Range("E7").Value = "  -0.37%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").Value = "'30.50"
$ws.Range("E10").Value = "  -1.16%  "

$ws.Range("E11").Value = "  -2.25%  "

$ws.Range("E12").Value = "  +3.38%  "

$ws.Range("E13").Value = "  +0.50%  "

$ws.Range("D14").Value = "2.571.01"
$ws.Range("E14").Value = "  -0.63%  "

$ws.Range("E15").Value = "  -1.94%  "

$ws.Range("D16").Value = "2.227.14"
$ws.Range("E16").Value = "  -1.15%  "

$ws.Range("E17").Value = "  -0.11%  "

$ws.Range("D18").Value = "40.084.29"
$ws.Range("E18").Value = "  +0.36%  "

$ws.Range("D19").Value = "0.0₃0890"
$ws.Range("E19").Value = "  -0.56%  "

$ws.Range("D20").Value = "'11.40"
$ws.Range("E20").Value = "  +7.99%  "

$ws.Range("E21").Value = "  -0.29%  "

$ws.Range("D22").Value = "'65.77"
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("D23").Value = "'237.39"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  +0.89%  "

$ws.Range("E26").Value = "  -0.61%  "

$ws.Range("D27").Value = "'22.80"
$ws.Range("E27").Value = "  -1.04%  "

$ws.Range("D28").Value = "'2.18"
$ws.Range("E28").Value = "  -1.63%  "

$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("D30").Value = "'156.08"
$ws.Range("E30").Value = "  +1.28%  "

$ws.Range("D31").Value = "'31.89"
$ws.Range("E31").Value = "  -6.44%  "

$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("D33").Value = "'4.96"
$ws.Range("E33").Value = "  +1.69%  "

$ws.Range("D34").Value = "'0.0721"
$ws.Range("E34").Value = "  +0.99%  "

$ws.Range("D35").Value = "'2.91"
$ws.Range("E35").Value = "  +7.30%  "

$ws.Range("D36").Value = "'2.35"
$ws.Range("E36").Value = "  -1.40%  "

$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").Value = "'15.82"
$ws.Range("E38").Value = "  -4.13%  "

$ws.Range("D39").Value = "'0.0987"
$ws.Range("E39").Value = "  -1.47%  "

$ws.Range("E40").Value = "  +1.76%  "

$ws.Range("D41").Value = "2.127.84"
$ws.Range("E41").Value = "  +8.00%  "

$ws.Range("D42").Value = "'3.89"
$ws.Range("E42").Value = "  +2.32%  "

$ws.Range("D43").Value = "'18.50"
$ws.Range("E43").Value = "  +13.71%  "

$ws.Range("E44").Value = "  -3.26%  "

$ws.Range("E45").Value = "  -0.96%  "

$ws.Range("D46").Value = "'9.87"
$ws.Range("E46").Value = "  +0.73%  "

$ws.Range("E47").Value = "  +3.78%  "

$ws.Range("D48").Value = "2.437.72"
$ws.Range("E48").Value = "  -0.84%  "

$ws.Range("E49").Value = "  +0.66%  "

$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").Value = "'1.11"
$ws.Range("E50").Value = "  +2.32%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'69.46"
$ws.Range("E51").Value = "  -2.26%  "
